$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 WHERE diag.primary_disease_site IN ['Lung']

OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp:sample)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
       coalesce(f.file_name, '') AS `File Name`,
       coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@

$ws.Range("B4").Value = $newFilesTabQuery
